$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1439.3726
$ws.Range("I137").Value = 1525.6
$ws.Range("J137").Value = 1316.1904
$ws.Range("K137").Value = 4576.799999999999
$ws.Range("L137").Value = 3948.5712
$ws.Range("M137").Value = -2026.799999999999
$ws.Range("N137").Value = -9048.5712

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1876.7142
$ws.Range("I2").Value = 1548.6
$ws.Range("J2").Value = 2697
$ws.Range("K2").Value = 1548.6
$ws.Range("L2").Value = 2697
$ws.Range("M2").Value = -1435.6
$ws.Range("N2").Value = -2923

$ws.Range("H45").Value = 14974.571
$ws.Range("I45").Value = 14974.571
$ws.Range("K45").Value = 14974.571
$ws.Range("M45").Value = -14597.571

$ws.Range("H61").Value = 5527.5483
$ws.Range("I61").Value = 6420.636
$ws.Range("J61").Value = 3344.4443
$ws.Range("K61").Value = 6420.636
$ws.Range("L61").Value = 3344.4443
$ws.Range("M61").Value = -6208.636
$ws.Range("N61").Value = -3768.4443

$ws.Range("H74").Value = 2021.8235
$ws.Range("I74").Value = 1955.909
$ws.Range("J74").Value = 2142.6667
$ws.Range("K74").Value = 1955.909
$ws.Range("L74").Value = 2142.6667
$ws.Range("M74").Value = -1081.909
$ws.Range("N74").Value = -3890.6667

$ws.Range("H77").Value = 2021.8235
$ws.Range("I77").Value = 1955.909
$ws.Range("J77").Value = 2142.6667
$ws.Range("K77").Value = 9779.545
$ws.Range("L77").Value = 10713.3335
$ws.Range("M77").Value = -5411.545
$ws.Range("N77").Value = -19449.3335

$ws.Range("H88").Value = 2142.5334
$ws.Range("I88").Value = 2060.125
$ws.Range("J88").Value = 2236.7144
$ws.Range("K88").Value = 2060.125
$ws.Range("L88").Value = 2236.7144
$ws.Range("M88").Value = -1654.125
$ws.Range("N88").Value = -3048.7144

$ws.Range("H91").Value = 2142.5334
$ws.Range("I91").Value = 2060.125
$ws.Range("J91").Value = 2236.7144
$ws.Range("K91").Value = 2060.125
$ws.Range("L91").Value = 2236.7144
$ws.Range("M91").Value = -656.125
$ws.Range("N91").Value = -5044.7144

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H116").Value = 1876.7142
$ws.Range("I116").Value = 1548.6
$ws.Range("J116").Value = 2697
$ws.Range("K116").Value = 1548.6
$ws.Range("L116").Value = 2697
$ws.Range("M116").Value = 745.4000000000001
$ws.Range("N116").Value = -7285

$ws.Range("H122").Value = 1069980.4
$ws.Range("I122").Value = 1509996.5
$ws.Range("J122").Value = 1370
$ws.Range("K122").Value = 4529989.5
$ws.Range("L122").Value = 4110
$ws.Range("M122").Value = -4527539.5
$ws.Range("N122").Value = -9010

$ws.Range("H136").Value = 5527.5483
$ws.Range("I136").Value = 6420.636
$ws.Range("J136").Value = 3344.4443
$ws.Range("K136").Value = 19261.908
$ws.Range("L136").Value = 10033.3329
$ws.Range("M136").Value = -16711.908
$ws.Range("N136").Value = -15133.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1876.7142
$ws.Range("I3").Value = 1548.6
$ws.Range("J3").Value = 2697
$ws.Range("K3").Value = 1548.6
$ws.Range("L3").Value = 2697
$ws.Range("M3").Value = -1434.6
$ws.Range("N3").Value = -2925

$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H86").Value = 1714
$ws.Range("I86").Value = 1674.1538
$ws.Range("J86").Value = 1800.3334
$ws.Range("K86").Value = 1674.1538
$ws.Range("L86").Value = 1800.3334
$ws.Range("M86").Value = -551.1538
$ws.Range("N86").Value = -4046.3334

$ws.Range("H89").Value = 1714
$ws.Range("I89").Value = 1674.1538
$ws.Range("J89").Value = 1800.3334
$ws.Range("K89").Value = 8370.769
$ws.Range("L89").Value = 9001.666999999999
$ws.Range("M89").Value = -2754.769
$ws.Range("N89").Value = -20233.667

$ws.Range("H134").Value = 6900.304
$ws.Range("I134").Value = 9407.929
$ws.Range("J134").Value = 2999.5557
$ws.Range("K134").Value = 28223.787
$ws.Range("L134").Value = 8998.667099999999
$ws.Range("M134").Value = -25688.787
$ws.Range("N134").Value = -14068.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 526825.75
$ws.Range("I113").Value = 540.75
$ws.Range("J113").Value = 909578.4399999999
$ws.Range("K113").Value = 1622.25
$ws.Range("L113").Value = 2728735.32
$ws.Range("M113").Value = 547.75
$ws.Range("N113").Value = -2733075.32

$ws.Range("H121").Value = 11731.4
$ws.Range("J121").Value = 28702
$ws.Range("L121").Value = 86106
$ws.Range("N121").Value = -88726

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3072.4443
$ws.Range("I80").Value = 2017.6666
$ws.Range("J80").Value = 3599.8333
$ws.Range("K80").Value = 2017.6666
$ws.Range("L80").Value = 3599.8333
$ws.Range("M80").Value = -1019.6666
$ws.Range("N80").Value = -5595.8333

$ws.Range("H83").Value = 3072.4443
$ws.Range("I83").Value = 2017.6666
$ws.Range("J83").Value = 3599.8333
$ws.Range("K83").Value = 10088.333
$ws.Range("L83").Value = 17999.1665
$ws.Range("M83").Value = -5096.333000000001
$ws.Range("N83").Value = -27983.1665

$ws.Range("H102").Value = 530807.5600000001
$ws.Range("J102").Value = 1525.5834
$ws.Range("L102").Value = 1525.5834
$ws.Range("N102").Value = -4769.5834

$ws.Range("H122").Value = 35443076
$ws.Range("I122").Value = 73414500
$ws.Range("J122").Value = 3079.2
$ws.Range("K122").Value = 220243500
$ws.Range("L122").Value = 9237.599999999999
$ws.Range("M122").Value = -220241050
$ws.Range("N122").Value = -14137.6

$ws.Range("H123").Value = 20600.793
$ws.Range("J123").Value = 20600.793
$ws.Range("L123").Value = 20600.793
$ws.Range("N123").Value = -25500.793

$ws.Range("H126").Value = 4872.6343
$ws.Range("I126").Value = 10855.637
$ws.Range("J126").Value = 2678.8667
$ws.Range("K126").Value = 32566.911
$ws.Range("L126").Value = 8036.6001
$ws.Range("M126").Value = -30096.911
$ws.Range("N126").Value = -12976.6001

$ws.Range("H132").Value = 2729
$ws.Range("I132").Value = 2487.625
$ws.Range("J132").Value = 2866.9285
$ws.Range("K132").Value = 7462.875
$ws.Range("L132").Value = 8600.7855
$ws.Range("M132").Value = -4932.875
$ws.Range("N132").Value = -13660.7855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 66669704
$ws.Range("I40").Value = 100001760
$ws.Range("K40").Value = 100001760
$ws.Range("M40").Value = -100001624

$ws.Range("H122").Value = 4288877.5
$ws.Range("I122").Value = 10206710
$ws.Range("J122").Value = 836809
$ws.Range("K122").Value = 30620130
$ws.Range("L122").Value = 2510427
$ws.Range("M122").Value = -30617680
$ws.Range("N122").Value = -2515327

$ws.Range("H132").Value = 22231864
$ws.Range("I132").Value = 55573810
$ws.Range("J132").Value = 3899.889
$ws.Range("K132").Value = 166721430
$ws.Range("L132").Value = 11699.667
$ws.Range("M132").Value = -166718900
$ws.Range("N132").Value = -16759.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1257.7142
$ws.Range("I122").Value = 1221
$ws.Range("J122").Value = 1306.6666
$ws.Range("K122").Value = 3663
$ws.Range("L122").Value = 3919.9998
$ws.Range("M122").Value = -1213
$ws.Range("N122").Value = -8819.9998

$ws.Range("H136").Value = 2431.8096
$ws.Range("I136").Value = 3080.4443
$ws.Range("J136").Value = 1945.3334
$ws.Range("K136").Value = 9241.332900000001
$ws.Range("L136").Value = 5836.0002
$ws.Range("M136").Value = -6691.332900000001
$ws.Range("N136").Value = -10936.0002
